$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the two rows that held the empty "[]" parameter sets
# (previously rows 26 and 27, labelled "TestCase10"); this shifts
# all subsequent rows up by two.
$ws.Rows("26:27").Delete()

# Relabel the remaining rows so the "Test Case Name" groups are
# contiguous and correctly numbered (TestCase07 / TestCase08).
$ws.Range("A20").Value = "TestCase07"
$ws.Range("A22").Value = "TestCase07"
$ws.Range("A23").Value = "TestCase08"
$ws.Range("A24").Value = "TestCase08"
$ws.Range("A25").Value = "TestCase08"

# Nudge the alignment so Excel materialises the "applyAlignment"
# cell style (matching the style already used by later groups) on
# the first/continuation rows of each relabelled group.
$ws.Range("A20").WrapText = $False
$ws.Range("A22").WrapText = $False
$ws.Range("A23").WrapText = $False
$ws.Range("A24").WrapText = $False
$ws.Range("A25").WrapText = $False
